$d = $word.ActiveDocument
$d.Content.Find.Execute(", SQL or ", $true, $false, $false, $false, $false, $true, 1, $false, ", SQL with ", 2)
